$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply all cell value changes described by the diff.
#
# All touched cells in this sheet hold plain text (OOXML t="inlineStr") values,
# even when their content looks like a number (e.g. "320.17", "1.000",
# "29.897.81"). If such a string is assigned straight to .Value, Excel's
# automatic type detection can silently reinterpret it as a numeric value,
# dropping significant trailing zeros and/or introducing floating point noise
# (e.g. "1.000" -> 1, "320.17" -> 320.17000000000002). To guarantee every cell
# keeps holding the exact original text, the whole target range is switched to
# the Text number format before any values are written, and the formatting is
# cleared again afterwards so the saved file does not keep any incidental
# style/number-format change on these cells.
$editRange = $ws.Range("D2:E51")
$editRange.NumberFormat = "@"

$ws.Range("D2").Value = '29.897.81'
$ws.Range("E2").Value = '  -1.24%  '
$ws.Range("D3").Value = '1.922.14'
$ws.Range("E3").Value = '  +1.53%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '320.17'
$ws.Range("E5").Value = '  -0.86%  '
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").Value = '0.5055'
$ws.Range("E7").Value = '  -2.36%  '
$ws.Range("D8").Value = '0.4051'
$ws.Range("E8").Value = '  +0.90%  '
$ws.Range("D9").Value = '0.08347'
$ws.Range("E9").Value = '  -0.83%  '
$ws.Range("D10").Value = '42.36'
$ws.Range("E10").Value = '  -0.79%  '
$ws.Range("D11").Value = '1.104'
$ws.Range("E11").Value = '  -0.90%  '
$ws.Range("D12").Value = '23.83'
$ws.Range("E12").Value = '  +3.32%  '
$ws.Range("D13").Value = '1.920.83'
$ws.Range("E13").Value = '  +1.84%  '
$ws.Range("D14").Value = '6.411'
$ws.Range("E14").Value = '  -0.33%  '
$ws.Range("D15").Value = '7.247'
$ws.Range("E15").Value = '  -0.90%  '
$ws.Range("D16").Value = '0.9967'
$ws.Range("E16").Value = '  -0.47%  '
$ws.Range("E17").Value = '  -2.06%  '
$ws.Range("D18").Value = '0.00001098'
$ws.Range("E18").Value = '  -0.96%  '
$ws.Range("D19").Value = '0.06517'
$ws.Range("E19").Value = '  -1.98%  '
$ws.Range("D20").Value = '18.26'
$ws.Range("E20").Value = '  +0.25%  '
$ws.Range("E21").Value = '  -0.17%  '
$ws.Range("D22").Value = '5.950'
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("D23").Value = '29.966.94'
$ws.Range("E23").Value = '  -0.97%  '
$ws.Range("D24").Value = '11.33'
$ws.Range("E24").Value = '  +0.39%  '
$ws.Range("E25").Value = '  -1.64%  '
$ws.Range("D26").Value = '2.160.40'
$ws.Range("E26").Value = '  +2.50%  '
$ws.Range("D27").Value = '22.08'
$ws.Range("E27").Value = '  +2.45%  '
$ws.Range("D28").Value = '162.26'
$ws.Range("E28").Value = '  +0.10%  '
$ws.Range("D29").Value = '2.337'
$ws.Range("E29").Value = '  +0.20%  '
$ws.Range("D30").Value = '128.89'
$ws.Range("E30").Value = '  -0.27%  '
$ws.Range("E31").Value = '  +4.16%  '
$ws.Range("E32").Value = '  -1.44%  '
$ws.Range("D33").Value = '5.957'
$ws.Range("E33").Value = '  -2.35%  '
$ws.Range("D34").Value = '3.773'
$ws.Range("E34").Value = '  +0.75%  '
$ws.Range("D35").Value = '5.419'
$ws.Range("E35").Value = '  +1.64%  '
$ws.Range("D36").Value = '0.02449'
$ws.Range("E36").Value = '  -1.82%  '
$ws.Range("D37").Value = '0.06420'
$ws.Range("E37").Value = '  -1.84%  '
$ws.Range("B38").Value = 'TheSandbox'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D38").Value = '0.6611'
$ws.Range("E38").Value = '  +1.66%  '
$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").Value = '0.2156'
$ws.Range("E39").Value = '  -1.86%  '
$ws.Range("D40").Value = '8.761'
$ws.Range("E40").Value = '  -0.55%  '
$ws.Range("D41").Value = '1.197'
$ws.Range("E41").Value = '  -1.90%  '
$ws.Range("D42").Value = '11.39'
$ws.Range("E42").Value = '  -3.47%  '
$ws.Range("D43").Value = '1.212'
$ws.Range("E43").Value = '  -1.13%  '
$ws.Range("D44").Value = '2.228'
$ws.Range("E44").Value = '  +8.56%  '
$ws.Range("D45").Value = '13.48'
$ws.Range("E45").Value = '  +2.24%  '
$ws.Range("D46").Value = '0.6105'
$ws.Range("E46").Value = '  +0.33%  '
$ws.Range("D47").Value = '3.614'
$ws.Range("E47").Value = '  -1.87%  '
$ws.Range("D48").Value = '1.210'
$ws.Range("E48").Value = '  -2.05%  '
$ws.Range("D49").Value = '121.80'
$ws.Range("E49").Value = '  -2.13%  '
$ws.Range("D50").Value = '79.06'
$ws.Range("E50").Value = '  +0.00%  '
$ws.Range("D51").Value = '1.128'
$ws.Range("E51").Value = '  -2.68%  '

$editRange.ClearFormats()
